$wb = $excel.ActiveWorkbook

# Target state for data rows 2-23 (columns B..I). Column A (the numeric index)
# is left untouched - it already holds the correct sequential values.
# Row 2 of the original sheet (the cancelled 2024-03-08 cosplay show) is
# removed by shifting every subsequent row's B..I content up by one and then
# deleting the now-superfluous last row (old row 24).
$rowsData = @(
    @(2, '2024-03-16', '苏州·OrangeOrange新春随舞派对【免费展会】', '狮山路298号 金鹰国际购物中心(狮山路店)', '2024.03.16 13:00-03.16 17:30', 133, 25, 'https://show.bilibili.com/platform/detail.html?id=82046', '//i2.hdslb.com/bfs/openplatform/202402/0OH3Ax4I1708913393518.png'),
    @(3, '2024-03-17', '苏州·世纪幻想动漫游戏展2.0', '清禾路886号 尹山湖大剧院', '2024.03.17 10:00-03.17 17:00', 1655, 60, 'https://show.bilibili.com/platform/detail.html?id=81387', '//i0.hdslb.com/bfs/openplatform/202402/isVyI9hH1708590817616.jpeg'),
    @(4, '2024-03-23', '苏州·Look Look动漫嘉年华', '阳澄半岛慈云路168号(重元寺北) 阳澄湖澜廷度假酒店', '2024.03.23 10:00-03.23 17:30', 639, 52.2, 'https://show.bilibili.com/platform/detail.html?id=81698', '//i1.hdslb.com/bfs/openplatform/202402/CP95X8ao1708934930351.jpeg'),
    @(5, '2024-03-30', '苏州·奇幻世界5.3动漫游戏展', '龙河路1288号 乐动力苏州湾体育中心', '2024.03.30 10:00-03.31 17:00', 1107, 55, 'https://show.bilibili.com/platform/detail.html?id=82002', '//i0.hdslb.com/bfs/openplatform/202402/HlxVHAz91708593664222.jpeg'),
    @(6, '2024-03-31', '张家港·META萌圆饿了', '大新镇人民路18号 新香苑宴会厅', '2024.03.31 10:00-03.31 17:00', 21, 30, 'https://show.bilibili.com/platform/detail.html?id=82407', '//i0.hdslb.com/bfs/openplatform/202403/pxZkSPcL1709707210340.jpeg'),
    @(7, '2024-04-04', '【会员购严选】苏州·Anime LIVE 国际动漫品牌博览会', '金山南路288号木渎影视城F2 苏州广电国际会展中心', '2024.04.04 10:00-04.05 17:00', 11642, 60, 'https://show.bilibili.com/platform/detail.html?id=81827', '//i2.hdslb.com/bfs/openplatform/202402/6oSFbWOx1707301464970.jpeg'),
    @(8, '2024-04-04', '常熟·SL动漫展', '报慈北路218号 四季花园酒店(报慈北路店)', '2024.04.04 10:00-04.04 17:00', 32, 50, 'https://show.bilibili.com/platform/detail.html?id=82250', '//i0.hdslb.com/bfs/openplatform/202403/rcvwgj7N1709518723294.jpeg'),
    @(9, '2024-04-06', '苏州·第一届寒假动漫展宅舞比赛-CF01', '润元路润南巷172号,地铁二号线陆慕站东200米,近市旅游换乘中心北100米 斐利酒店', '2024.04.06 10:00-04.06 16:00', 94, 49, 'https://show.bilibili.com/platform/detail.html?id=80528', '//i2.hdslb.com/bfs/openplatform/202401/oWbVnOjD1704445446390.jpeg'),
    @(10, '2024-04-13', '苏州·X-party 国漫游戏嘉年华03', '秋枫街与开平路交叉口西南角 爱琴海购物中心', '2024.04.13 10:00-04.14 17:00', 462, 48, 'https://show.bilibili.com/platform/detail.html?id=82042', '//i2.hdslb.com/bfs/openplatform/202403/GWNvc78z1709275224442.jpeg'),
    @(11, '2024-04-13', '苏州·绘时国乙1.0-秩序之外', '石路步行街永福桥浜15号 银河广场', '2024.04.13 13:30-04.13 20:00', 381, 88, 'https://show.bilibili.com/platform/detail.html?id=80789', '//i0.hdslb.com/bfs/openplatform/202403/nIPoXWqO1709275656198.jpeg'),
    @(12, '2024-04-20', '苏州·首届Redamancy动漫游戏嘉年华', '清禾路886号 尹山湖大剧院', '2024.04.20 10:00-04.20 17:00', 1098, 60, 'https://show.bilibili.com/platform/detail.html?id=81879', '//i0.hdslb.com/bfs/openplatform/202402/lR4oJWzI1708309129629.jpeg'),
    @(13, '2024-04-21', '苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0', '清禾路888号2号楼3楼 格莱美婚礼宴会中心', '2024.04.21 10:00-04.21 21:00', 819, 69.9, 'https://show.bilibili.com/platform/detail.html?id=78666', '//i0.hdslb.com/bfs/openplatform/202312/X0PZ3YhH1703822037665.jpeg'),
    @(14, '2024-05-01', '昆山·第十二届理想乡动漫游戏展', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.01 10:00-05.03 17:00', 12990, 75, 'https://show.bilibili.com/platform/detail.html?id=77196', '//i2.hdslb.com/bfs/openplatform/202310/9xMTQMlg1696736126094.png'),
    @(15, '2024-05-01', '苏州·第十七届 I COME ACG  动漫品牌博览会', '金山南路288号 广电国际会展中心', '2024.05.01 10:00-05.02 17:00', 13206, 65, 'https://show.bilibili.com/platform/detail.html?id=79789', '//i2.hdslb.com/bfs/openplatform/202312/lau3mW031702535438289.jpeg'),
    @(16, '2024-05-02', '昆山·第十二届理想乡动漫游戏展嘉宾北齐后主签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.02 14:00-05.02 16:00', 35, 1, 'https://show.bilibili.com/platform/detail.html?id=81116', '//i2.hdslb.com/bfs/openplatform/202401/EubrAneC1705648695005.jpeg'),
    @(17, '2024-05-02', '昆山·第十二届理想乡动漫游戏展嘉宾啊川签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.02 14:00-05.02 16:00', 146, 1, 'https://show.bilibili.com/platform/detail.html?id=81100', '//i2.hdslb.com/bfs/openplatform/202401/F24i5GMX1705646667852.jpeg'),
    @(18, '2024-05-02', '昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.02 14:00-05.02 16:00', 18, 1, 'https://show.bilibili.com/platform/detail.html?id=81119', '//i2.hdslb.com/bfs/openplatform/202401/SDnLB1gR1705648838683.jpeg'),
    @(19, '2024-05-02', '昆山·第十二届理想乡动漫游戏展嘉宾葫芦岛老八签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.02 14:00-05.02 16:00', 36, 1, 'https://show.bilibili.com/platform/detail.html?id=81118', '//i2.hdslb.com/bfs/openplatform/202401/uHOCneLv1705648779163.jpeg'),
    @(20, '2024-05-03', '常熟·CDW·动漫展03', '常熟国际展览中心 国际展览中心', '2024.05.03 09:00-05.04 17:30', 257, 60, 'https://show.bilibili.com/platform/detail.html?id=82489', '//i0.hdslb.com/bfs/openplatform/202403/XK411blC1709794808211.jpeg'),
    @(21, '2024-05-03', '昆山·第十二届理想乡动漫游戏展嘉宾沈辞签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.03 14:00-05.03 16:00', 89, 1, 'https://show.bilibili.com/platform/detail.html?id=81120', '//i0.hdslb.com/bfs/openplatform/202401/4Pay1rR61705648901961.jpeg'),
    @(22, '2024-05-03', '昆山·第十二届理想乡动漫游戏展嘉宾矮乐多aliga签售会', '花桥经济开发区绿地大道1598号 花桥国际博览中心', '2024.05.03 14:00-05.03 16:00', 45, 1, 'https://show.bilibili.com/platform/detail.html?id=81114', '//i1.hdslb.com/bfs/openplatform/202401/Peub7FOc1705648580577.jpeg'),
    @(23, '2024-06-08', '【会员购严选】苏州·Come in joy动漫国潮文化节', '金山南路288号 广电国际会展中心', '2024.06.08 10:00-06.09 17:00', 136, 60, 'https://show.bilibili.com/platform/detail.html?id=82233', '//i0.hdslb.com/bfs/openplatform/202403/F86lgbSt1709278264141.jpeg')
)

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($entry in $rowsData) {
        $r = $entry[0]

        # Column B holds a plain "YYYY-MM-DD" text label. Excel's COM layer
        # auto-detects that shape as a date on assignment, so force the cell
        # to Text first and strip the format back off afterwards (ClearFormats
        # removes the now-unneeded "@" number format but keeps the text type
        # and value that was already committed).
        $bCell = $ws.Cells.Item($r, 2)
        $bCell.NumberFormat = "@"
        $bCell.Value = $entry[1]
        $bCell.ClearFormats()

        $ws.Cells.Item($r, 3).Value = $entry[2]
        $ws.Cells.Item($r, 4).Value = $entry[3]
        $ws.Cells.Item($r, 5).Value = $entry[4]
        $ws.Cells.Item($r, 6).Value = $entry[5]
        $ws.Cells.Item($r, 7).Value = $entry[6]
        $ws.Cells.Item($r, 8).Value = $entry[7]
        $ws.Cells.Item($r, 9).Value = $entry[8]
    }

    # The last remaining row (old row 24) has now been duplicated into row 23;
    # remove the original row 24 entirely so the sheet ends at row 23.
    $ws.Rows.Item(24).Delete()
}
